$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.892.07"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").Value = "3.556.25"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.44"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "647.96"
$ws.Range("E6").Value = "  +1.59%  "

$ws.Range("E7").Value = "  -1.27%  "

$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -2.30%  "

$ws.Range("D11").Value = "3.556.12"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.19"
$ws.Range("E13").Value = "  -2.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.47"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "4.216.58"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").Value = "94.898.72"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").Value = "3.548.50"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("E20").Value = "  -4.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.82"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.44"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "504.71"
$ws.Range("E23").Value = "  -2.23%  "

$ws.Range("E24").Value = "  -5.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.76"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.07"
$ws.Range("E27").Value = "  -1.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.41"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "3.747.29"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.35"
$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.177"
$ws.Range("E35").Value = "  -3.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.61"
$ws.Range("E36").Value = "  +4.62%  "

$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("E38").Value = "  +7.08%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.47"
$ws.Range("E39").Value = "  +6.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "579.18"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  -1.27%  "

$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.76"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("E45").Value = "  +4.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.64"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.28"
$ws.Range("E48").Value = "  +31.64%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0411"
$ws.Range("E49").Value = "  -5.40%  "

$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.23"
$ws.Range("E51").Value = "  -1.37%  "
